# Updated symbol list on Sun Dec 18 15:21:47 UTC 2022 with GitHub Actions
# Applies the coin price / listing refresh described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $address, $value) {
    # Values in column D look numeric (e.g. "246.79") but must be stored as
    # literal text, exactly as the source sheet already has them (inline
    # strings). Forcing the cell to Text format before assignment prevents
    # Excel from re-interpreting the string as a floating point number
    # (which would introduce binary rounding noise such as 246.78999999999999).
    $cell = $sheet.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---- Column D (Price) updates for otherwise-unchanged rows ----
Set-TextValue $ws "D2"  "246.79"
Set-TextValue $ws "D3"  "22.42"
Set-TextValue $ws "D4"  "5.462"
Set-TextValue $ws "D5"  "0.05637"
Set-TextValue $ws "D7"  "0.8053"
Set-TextValue $ws "D9"  "0.1435"
Set-TextValue $ws "D10" "0.07339"
Set-TextValue $ws "D11" "0.03189"
Set-TextValue $ws "D12" "0.02933"
Set-TextValue $ws "D13" "0.09264"
Set-TextValue $ws "D14" "0.001664"
Set-TextValue $ws "D15" "3.210"
Set-TextValue $ws "D16" "0.04726"
Set-TextValue $ws "D17" "0.0005829"
Set-TextValue $ws "D18" "0.006472"
Set-TextValue $ws "D19" "0.001055"
Set-TextValue $ws "D20" "0.004108"
Set-TextValue $ws "D21" "0.0001504"
Set-TextValue $ws "D22" "3.982"
Set-TextValue $ws "D23" "3.387"
Set-TextValue $ws "D24" "2.142"
Set-TextValue $ws "D25" "0.3267"
Set-TextValue $ws "D26" "0.1318"
Set-TextValue $ws "D27" "0.0002913"
Set-TextValue $ws "D40" "0.04151"
Set-TextValue $ws "D44" "0.009035"
Set-TextValue $ws "D45" "0.00005644"
Set-TextValue $ws "D46" "0.00000000753"
Set-TextValue $ws "D47" "0.6823"

# ---- Rows 41-43: the three coins (KickToken / BKEXToken / CEJI) cycled
#      positions (row 42 -> 41, row 43 -> 42, row 41 -> 43) and picked up
#      refreshed Price / rank-label values. ----

# Row 41 becomes what used to be BKEXToken (row 42), with a new price.
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D41" "0.1039"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 becomes what used to be CEJI (row 43), with a new price.
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.002978"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 becomes what used to be KickToken (row 41), with a new price, and
# is now flagged as the day's worst performer.
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D43" "0.003247"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# ---- Row 48 (BOLO): refreshed price and loses the "worst performer" flag ----
Set-TextValue $ws "D48" "0.01858"
$ws.Range("E48").Value = "47BOLOBOLO"

Write-Host "Applied coin price/listing refresh"
